# The workbook "Saldo" lists account balances. This edit removes a set of
# rows (accounts) from the "Export" sheet's data table, shifting the
# remaining rows up, exactly matching the rows removed in the source diff:
#
#   row 2  -> 004213929 RODOLFO    300000
#   row 3  -> 004474776 GILSON     123642.88
#   row 4  -> 004940699 RACHEL     58909.87
#   row 5  -> 004936634 LEONARDO   22035.04
#   row 6  -> 005046805 RICARDO    7000
#   row 8  -> 004479965 DIEGO      6102.06
#   row 9  -> 004420763 CHRISTIAN  5483.53
#   row 11 -> 002823185 SIMONE     3978.2
#   row 14 -> 005654767 DIEGO      1128.97
#
# Row 1 is the header ("Conta", "Nome", "Saldo") and is left untouched.
# Deletions are applied from the bottom of the sheet upward so that row
# numbers of not-yet-processed rows remain stable while earlier rows are
# removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14:A14").EntireRow.Delete()
$ws.Range("A11:A11").EntireRow.Delete()
$ws.Range("A8:A9").EntireRow.Delete()
$ws.Range("A2:A6").EntireRow.Delete()
